# Apply localization sheet updates:
#  - Normalize placeholder text from "{0}" to "XXXX" for the untranslated
#    language columns (C:K) on the LOGIN_ERROR, PAGE, and REGISTER_ERROR rows
#  - Update the frozen-pane scroll position and active selection on the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose C:K cells still contained the literal "{0}" placeholder;
# replace with the standard "XXXX" untranslated marker used elsewhere.
$rows = @(20, 29, 34)
foreach ($r in $rows) {
    $ws.Range("C$r`:K$r").Value = "XXXX"
}

# Update the view state: scroll/freeze pane top-left cell and the active selection.
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("C19:C35").Select()
